$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header columns: DefaultName -> Name, DefaultPrice -> Price, DefaultExplain -> Explain
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Price"
$ws.Range("E1").Value = "Explain"
